$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 columns for apoio_std/apoio_min/apoio_max right after apoio_medio (col L)
$ws.Range("M1:O1").EntireColumn.Insert()
# Insert 3 columns for contribuicoes_std/min/max right after contribuicoes_med (now col Q)
$ws.Range("R1:T1").EntireColumn.Insert()

# Rename headers that shifted semantics in place (values unchanged, just renamed)
$ws.Range("H1").Value = "arrecadado_avg"
$ws.Range("I1").Value = "arrecadado_std"
$ws.Range("J1").Value = "arrecadado_min"
$ws.Range("K1").Value = "arrecadado_max"
$ws.Range("M1").Value = "apoio_std"
$ws.Range("N1").Value = "apoio_min"
$ws.Range("O1").Value = "apoio_max"
$ws.Range("Q1").Value = "contribuicoes_med"
$ws.Range("R1").Value = "contribuicoes_std"
$ws.Range("S1").Value = "contribuicoes_min"
$ws.Range("T1").Value = "contribuicoes_max"

# Recalculated apoio_medio + new apoio_std/min/max + new contribuicoes_std/min/max per row
$ws.Cells.Item(2, 12).Value = 76.06649705332148
$ws.Cells.Item(2, 13).Value = 30.35289755311455
$ws.Cells.Item(2, 14).Value = 30.69891629110986
$ws.Cells.Item(2, 15).Value = 216.9152091005833
$ws.Cells.Item(2, 18).Value = 410.5455139511234
$ws.Cells.Item(2, 19).Value = 14
$ws.Cells.Item(2, 20).Value = 3474
$ws.Cells.Item(3, 12).Value = 82.12933179093685
$ws.Cells.Item(3, 13).Value = 35.40469491555042
$ws.Cells.Item(3, 14).Value = 33.97203328414528
$ws.Cells.Item(3, 15).Value = 230.5829516876651
$ws.Cells.Item(3, 18).Value = 174.8393084425899
$ws.Cells.Item(3, 19).Value = 6
$ws.Cells.Item(3, 20).Value = 1540
$ws.Cells.Item(4, 12).Value = 84.02563974325884
$ws.Cells.Item(4, 13).Value = 39.83943189124467
$ws.Cells.Item(4, 14).Value = 11.52676430516467
$ws.Cells.Item(4, 15).Value = 254.2443749773306
$ws.Cells.Item(4, 18).Value = 155.4234048603088
$ws.Cells.Item(4, 19).Value = 1
$ws.Cells.Item(4, 20).Value = 1330
$ws.Cells.Item(5, 12).Value = 83.58617223970138
$ws.Cells.Item(5, 13).Value = 40.20112754081283
$ws.Cells.Item(5, 14).Value = 16.05922078302265
$ws.Cells.Item(5, 15).Value = 386.987447085286
$ws.Cells.Item(5, 18).Value = 248.1423584049538
$ws.Cells.Item(5, 19).Value = 2
$ws.Cells.Item(5, 20).Value = 2120
$ws.Cells.Item(6, 12).Value = 83.40563308192627
$ws.Cells.Item(6, 13).Value = 32.56399533953332
$ws.Cells.Item(6, 14).Value = 18.47818326605706
$ws.Cells.Item(6, 15).Value = 195.20880750356
$ws.Cells.Item(6, 18).Value = 501.9791222296157
$ws.Cells.Item(6, 19).Value = 1
$ws.Cells.Item(6, 20).Value = 7954
$ws.Cells.Item(7, 12).Value = 77.97569054482099
$ws.Cells.Item(7, 13).Value = 35.65744130420693
$ws.Cells.Item(7, 14).Value = 10.77163914429046
$ws.Cells.Item(7, 15).Value = 195.6882025465182
$ws.Cells.Item(7, 18).Value = 459.3640728169867
$ws.Cells.Item(7, 19).Value = 1
$ws.Cells.Item(7, 20).Value = 7954
$ws.Cells.Item(8, 12).Value = 77.73968420752422
$ws.Cells.Item(8, 13).Value = 41.30988410144955
$ws.Cells.Item(8, 14).Value = 12.19662302883409
$ws.Cells.Item(8, 15).Value = 247.2901437851162
$ws.Cells.Item(8, 18).Value = 685.9026443808731
$ws.Cells.Item(8, 19).Value = 2
$ws.Cells.Item(8, 20).Value = 7954
$ws.Cells.Item(9, 12).Value = 81.75853347173708
$ws.Cells.Item(9, 13).Value = 38.14520603080047
$ws.Cells.Item(9, 14).Value = 18.47818326605706
$ws.Cells.Item(9, 15).Value = 234.707661751482
$ws.Cells.Item(9, 18).Value = 508.76757749687
$ws.Cells.Item(9, 19).Value = 2
$ws.Cells.Item(9, 20).Value = 7954
$ws.Cells.Item(10, 12).Value = 82.92407682444032
$ws.Cells.Item(10, 13).Value = 38.77988334228132
$ws.Cells.Item(10, 14).Value = 11.93343625774652
$ws.Cells.Item(10, 15).Value = 230.5829516876651
$ws.Cells.Item(10, 18).Value = 247.246130032236
$ws.Cells.Item(10, 19).Value = 1
$ws.Cells.Item(10, 20).Value = 2684
$ws.Cells.Item(11, 12).Value = 71.78666858221021
$ws.Cells.Item(11, 13).Value = 29.97870848948209
$ws.Cells.Item(11, 14).Value = 16.18065842403185
$ws.Cells.Item(11, 15).Value = 216.9152091005833
$ws.Cells.Item(11, 18).Value = 331.1306730481258
$ws.Cells.Item(11, 19).Value = 3
$ws.Cells.Item(11, 20).Value = 3474
$ws.Cells.Item(12, 12).Value = 71.78666858221021
$ws.Cells.Item(12, 13).Value = 29.97870848948209
$ws.Cells.Item(12, 14).Value = 16.18065842403185
$ws.Cells.Item(12, 15).Value = 216.9152091005833
$ws.Cells.Item(12, 18).Value = 331.1306730481258
$ws.Cells.Item(12, 19).Value = 3
$ws.Cells.Item(12, 20).Value = 3474
$ws.Cells.Item(13, 12).Value = 79.12786981308152
$ws.Cells.Item(13, 13).Value = 35.35639160943987
$ws.Cells.Item(13, 14).Value = 20.51363271354002
$ws.Cells.Item(13, 15).Value = 233.3973531230909
$ws.Cells.Item(13, 18).Value = 328.5928536530323
$ws.Cells.Item(13, 19).Value = 1
$ws.Cells.Item(13, 20).Value = 4584
$ws.Cells.Item(14, 12).Value = 84.32898346466456
$ws.Cells.Item(14, 13).Value = 39.18079463334893
$ws.Cells.Item(14, 14).Value = 11.52676430516467
$ws.Cells.Item(14, 15).Value = 195.6882025465182
$ws.Cells.Item(14, 18).Value = 147.3429131999576
$ws.Cells.Item(14, 19).Value = 1
$ws.Cells.Item(14, 20).Value = 612
$ws.Cells.Item(15, 12).Value = 83.48264574282582
$ws.Cells.Item(15, 13).Value = 40.20970245451376
$ws.Cells.Item(15, 14).Value = 21.00493274015408
$ws.Cells.Item(15, 15).Value = 247.2901437851162
$ws.Cells.Item(15, 18).Value = 174.0035582476068
$ws.Cells.Item(15, 19).Value = 2
$ws.Cells.Item(15, 20).Value = 770
$ws.Cells.Item(16, 12).Value = 83.64941179158359
$ws.Cells.Item(16, 13).Value = 36.99022540033587
$ws.Cells.Item(16, 14).Value = 16.18065842403185
$ws.Cells.Item(16, 15).Value = 254.2443749773306
$ws.Cells.Item(16, 18).Value = 183.1544319258093
$ws.Cells.Item(16, 19).Value = 1
$ws.Cells.Item(16, 20).Value = 1540
$ws.Cells.Item(17, 12).Value = 70.93306185876429
$ws.Cells.Item(17, 13).Value = 30.19803349932243
$ws.Cells.Item(17, 14).Value = 20.33774597757668
$ws.Cells.Item(17, 15).Value = 159.7763429092917
$ws.Cells.Item(17, 18).Value = 975.3935739169402
$ws.Cells.Item(17, 19).Value = 3
$ws.Cells.Item(17, 20).Value = 7954
$ws.Cells.Item(18, 12).Value = 76.84102373029619
$ws.Cells.Item(18, 13).Value = 33.37177139781743
$ws.Cells.Item(18, 14).Value = 16.18065842403185
$ws.Cells.Item(18, 15).Value = 226.5579622472015
$ws.Cells.Item(18, 18).Value = 457.1568742729124
$ws.Cells.Item(18, 19).Value = 1
$ws.Cells.Item(18, 20).Value = 7954
$ws.Cells.Item(19, 12).Value = 79.19230719197579
$ws.Cells.Item(19, 13).Value = 27.01736191709247
$ws.Cells.Item(19, 14).Value = 40.63189862969614
$ws.Cells.Item(19, 15).Value = 130.9739254174068
$ws.Cells.Item(19, 18).Value = 143.8139060050855
$ws.Cells.Item(19, 19).Value = 2
$ws.Cells.Item(19, 20).Value = 467
$ws.Cells.Item(20, 12).Value = 78.67160937524555
$ws.Cells.Item(20, 13).Value = 39.22105185666557
$ws.Cells.Item(20, 14).Value = 14.90596347946683
$ws.Cells.Item(20, 15).Value = 461.5197709071476
$ws.Cells.Item(20, 18).Value = 402.3874992420548
$ws.Cells.Item(20, 19).Value = 1
$ws.Cells.Item(20, 20).Value = 7954
$ws.Cells.Item(21, 12).Value = 69.51944033042635
$ws.Cells.Item(21, 13).Value = 28.22546328675653
$ws.Cells.Item(21, 14).Value = 17.82064921105857
$ws.Cells.Item(21, 15).Value = 196.4212117364618
$ws.Cells.Item(21, 18).Value = 448.9315388299796
$ws.Cells.Item(21, 19).Value = 4
$ws.Cells.Item(21, 20).Value = 4584
$ws.Cells.Item(22, 12).Value = 74.35971583315494
$ws.Cells.Item(22, 13).Value = 27.7214279498591
$ws.Cells.Item(22, 14).Value = 17.83984513748501
$ws.Cells.Item(22, 15).Value = 156.0426904908593
$ws.Cells.Item(22, 18).Value = 190.7586913141665
$ws.Cells.Item(22, 19).Value = 1
$ws.Cells.Item(22, 20).Value = 1879